$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B", "E", "F", "G")

$groups = @(
    ,@(142, 143)
    ,@(154, 155, 156)
    ,@(305, 306)
    ,@(343, 344)
    ,@(347, 348)
    ,@(364, 365)
    ,@(367, 368)
    ,@(371, 372)
    ,@(374, 375)
    ,@(392, 393)
    ,@(413, 414)
    ,@(578, 579)
    ,@(582, 583)
    ,@(585, 586)
    ,@(591, 592)
    ,@(679, 680)
)

foreach ($rows in $groups) {
    $n = $rows.Count
    $old = @{}
    foreach ($r in $rows) {
        $old[$r] = @{}
        foreach ($c in $cols) {
            $old[$r][$c] = $ws.Range("$c$r").Value2
        }
    }
    for ($i = 0; $i -lt $n; $i++) {
        $destRow = $rows[$i]
        $srcIndex = ($i - 1 + $n) % $n
        $srcRow = $rows[$srcIndex]
        foreach ($c in $cols) {
            $ws.Range("$c$destRow").Value = $old[$srcRow][$c]
        }
    }
}

Write-Host "Done rotating $($groups.Count) groups"